# "updates fixes most synchronization issues and works with rapid start"
#
# The bug/task tracker on Sheet1 is rebuilt: several rows are re-ordered,
# a handful of already-tracked bugs are re-prioritized (column D), and a
# big batch of new tasks (sync/cache/UI/rapid-start related) is appended
# below the existing list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old B3:D19 task table first -- the new layout reuses some of
# those row numbers for different text/priorities and leaves rows 15 and
# 18 blank, so starting from a clean slate avoids stray leftovers.
$ws.Range("B3:D40").ClearContents()

# --- existing tasks, re-ordered / re-prioritized -----------------------
$ws.Range("B3").Value  = 'Fix dumy tokens '
$ws.Range("C3").Value  = 'x'

$ws.Range("B4").Value  = 'fix Configuration menu'
$ws.Range("C4").Value  = 'x'
$ws.Range("D4").Value  = 2

$ws.Range("B5").Value  = 'fix finish message after user is done with setting up accounts'
$ws.Range("C5").Value  = 'x'

$ws.Range("B6").Value  = 'fix refresh of data from Dash services'
$ws.Range("C6").Value  = 'x'
$ws.Range("D6").Value  = 3

$ws.Range("B7").Value  = 'Create Logout button'
$ws.Range("C7").Value  = 'x'
$ws.Range("D7").Value  = 1

$ws.Range("B8").Value  = 'Create Format for saving User profile locally'
$ws.Range("C8").Value  = 'x'
$ws.Range("D8").Value  = 4

$ws.Range("B9").Value  = 'Create Parsing for local parse file'
$ws.Range("C9").Value  = 'x'
$ws.Range("D9").Value  = 5

$ws.Range("B10").Value = 'Make Sure all token submission processesare processed without bugs'
$ws.Range("C10").Value = 'x'

$ws.Range("B11").Value = 'Create reset functionality if user cancels out of Account bindign process'
$ws.Range("C11").Value = 'x'

$ws.Range("B12").Value = 'Write Php script to clear entry for ease'
$ws.Range("C12").Value = 'x'

$ws.Range("B13").Value = 'write php script to reset everything with app'
$ws.Range("C13").Value = 'x'

$ws.Range("B14").Value = 'back button without setting css to ui-dark.css'
$ws.Range("C14").Value = 'x'

# row 15 intentionally left blank

$ws.Range("B16").Value = 'user accounts'

$ws.Range("B17").Value = 'fix background notification bug caused by launcing application w/o internet and then logging in with internet'

# row 18 intentionally left blank

# --- new tasks added below the old list ---------------------------------
$ws.Range("B19").Value = 'Remove entry in database cache when service is unregistered.'
$ws.Range("D19").Value = 6

$ws.Range("B20").Value = 'Remove entry in cache when service is unregistered.'
$ws.Range("C20").Value = 'x'
$ws.Range("D20").Value = 4

$ws.Range("B21").Value = 'Remove entry in UI when service is unregistered.'
$ws.Range("C21").Value = 'x'
$ws.Range("D21").Value = 5

$ws.Range("B22").Value = 'Order Element in UI based on Time independent of the service'
$ws.Range("C22").Value = 'x'
$ws.Range("D22").Value = 7

$ws.Range("B23").Value = 'setup cron job for PHP'
$ws.Range("C23").Value = 'x'
$ws.Range("D23").Value = 2

$ws.Range("B24").Value = 'Investigate how you can have intelligent caching on server'
$ws.Range("C24").Value = 'x'
$ws.Range("D24").Value = 1

$ws.Range("B25").Value = 'See if you can use gmail atom to retrieve email'
$ws.Range("D25").Value = 8

$ws.Range("B26").Value = 'Setup groupon for cache and make appropriate calls'
$ws.Range("C26").Value = 'x'
$ws.Range("D26").Value = 3

$ws.Range("B27").Value = 'Fix look of setup when I right click while on settings menu.'
$ws.Range("D27").Value = 9

$ws.Range("B28").Value = 'Provide option for removing/adding image from cached rapid start storage'
$ws.Range("D28").Value = 10

$ws.Range("B29").Value = 'Fix default screen UI'
$ws.Range("D29").Value = 11

$ws.Range("B30").Value = 'Add other service like instagram yahoo and other deals'
$ws.Range("D30").Value = 14

$ws.Range("B31").Value = 'Fix dummy email'
$ws.Range("D31").Value = 12

$ws.Range("B32").Value = 'set loading bar for loading situations'
$ws.Range("D32").Value = 13

$ws.Range("B33").Value = 'fix getData function to have option for saving refreshed data to dbase'
$ws.Range("D33").Value = 14

$ws.Range("B34").Value = 'Provide ability to reply to tweets, facebookpost and email'
$ws.Range("D34").Value = 15

$ws.Range("B35").Value = 'Fix weird twitter bug of double login failure crash'
$ws.Range("C35").Value = 'x'
$ws.Range("D35").Value = 16

$ws.Range("B36").Value = 'Add Weather widget'
$ws.Range("D36").Value = 17

# --- view state: leave the selection where the author left off ----------
$ws.Range("C35").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
